$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Delete the "practiceQ" sheet (reported bug data no longer needed)
$practiceSheet = $wb.Worksheets.Item("practiceQ")
$practiceSheet.Delete()

# Rename remaining sheets to title case
$textEditorSheet = $wb.Worksheets.Item("textEditor")
$textEditorSheet.Name = "TextEditor"

$credentialsSheet = $wb.Worksheets.Item("credentials")
$credentialsSheet.Name = "Credentials"

# Make TextEditor the active/selected sheet with a new selection
$textEditorSheet.Activate()
$textEditorSheet.Range("B19").Select()

$excel.DisplayAlerts = $true
